$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Andalusia"
$ws.Range("B6").Value = "Balearics"
$ws.Range("B7").Value = "Canary Islands"
$ws.Range("B11").Value = "Catalonia"
$ws.Range("B12").Value = "Valencia"
$ws.Range("B17").Value = "Navarre"
$ws.Range("B18").Value = "Basque Country"

$ws.Range("B11").Select()
